# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay plain text
# (e.g. "28.345.29", "0.9964"), so force text formatting on the whole
# price column before writing, mirroring the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.345.29"
$ws.Range("E2").Value = "  -0.94%  "

# Row 3
$ws.Range("D3").Value = "1.825.88"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4
$ws.Range("D4").Value = "0.9964"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("D5").Value = "317.47"
$ws.Range("E5").Value = "  -0.17%  "

# Row 6
$ws.Range("D6").Value = "0.9964"
$ws.Range("E6").Value = "  -0.46%  "

# Row 7
$ws.Range("D7").Value = "0.5339"
$ws.Range("E7").Value = "  -1.83%  "

# Row 8
$ws.Range("D8").Value = "0.3993"
$ws.Range("E8").Value = "  +5.15%  "

# Row 9
$ws.Range("D9").Value = "0.07561"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").Value = "41.73"
$ws.Range("E10").Value = "  -1.17%  "

# Row 11
$ws.Range("D11").Value = "1.100"
$ws.Range("E11").Value = "  -1.21%  "

# Row 12
$ws.Range("D12").Value = "7.621"
$ws.Range("E12").Value = "  +3.73%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.273"
$ws.Range("E13").Value = "  +1.76%  "

# Row 14
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").Value = "0.9964"
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").Value = "20.65"
$ws.Range("E15").Value = "  -0.24%  "

# Row 16
$ws.Range("D16").Value = "1.811.98"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").Value = "89.56"
$ws.Range("E17").Value = "  -0.73%  "

# Row 18
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").Value = "0.06588"
$ws.Range("E19").Value = "  +0.56%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.40"
$ws.Range("E20").Value = "  -0.39%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.9971"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("D22").Value = "6.036"
$ws.Range("E22").Value = "  +1.34%  "

# Row 23
$ws.Range("D23").Value = "28.357.96"
$ws.Range("E23").Value = "  -0.94%  "

# Row 24
$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25
$ws.Range("D25").Value = "2.088"
$ws.Range("E25").Value = "  +0.61%  "

# Row 26
$ws.Range("D26").Value = "156.19"
$ws.Range("E26").Value = "  -3.39%  "

# Row 27
$ws.Range("D27").Value = "20.52"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.393"
$ws.Range("E28").Value = "  +2.25%  "

# Row 29
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.020.48"
$ws.Range("E29").Value = "  +0.67%  "

# Row 30
$ws.Range("D30").Value = "123.85"
$ws.Range("E30").Value = "  +0.69%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.109"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1097"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.07426"
$ws.Range("E33").Value = "  +12.43%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.669"
$ws.Range("E34").Value = "  -0.32%  "

# Row 35
$ws.Range("D35").Value = "5.585"
$ws.Range("E35").Value = "  -1.11%  "

# Row 36
$ws.Range("D36").Value = "0.2232"
$ws.Range("E36").Value = "  -1.25%  "

# Row 37
$ws.Range("D37").Value = "5.218"
$ws.Range("E37").Value = "  +3.32%  "

# Row 38
$ws.Range("D38").Value = "0.02302"
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("D39").Value = "8.647"
$ws.Range("E39").Value = "  +0.29%  "

# Row 40
$ws.Range("D40").Value = "11.33"
$ws.Range("E40").Value = "  +0.73%  "

# Row 41
$ws.Range("D41").Value = "0.6223"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("E42").Value = "  -0.46%  "

# Row 43
$ws.Range("E43").Value = "  -3.26%  "

# Row 44
$ws.Range("D44").Value = "13.41"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "3.687"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46
$ws.Range("D46").Value = "0.5795"
$ws.Range("E46").Value = "  -1.02%  "

# Row 47
$ws.Range("D47").Value = "125.18"
$ws.Range("E47").Value = "  -1.81%  "

# Row 48
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49
$ws.Range("D49").Value = "1.189"
$ws.Range("E49").Value = "  -0.02%  "

# Row 50
$ws.Range("D50").Value = "0.06867"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$ws.Range("D51").Value = "71.32"
$ws.Range("E51").Value = "  -1.54%  "
